# Generate Report for Handoff
# Adds two new localization files to the status report:
#   547613be-f236-4205-81d9-d6225ab7b667.md  (Ready for handoff)
#   d88e27aa-c46e-4436-a424-97a4e33712ca.md  (Ready for handoff)
# across the Overview, zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

$mdBase    = "https://github.com/OpenLocalizationTest/oltest/blob/40cadb3034932f4de1c0eca1b07fa12248308706/e2e/"
$zhcnBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1c1d52c119c85f2598fb34508664505d639282ad/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/"
$dedeBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/83dfd44bb66c316e4225af2c63d95b628ba523c2/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/"

$file1 = "547613be-f236-4205-81d9-d6225ab7b667"
$file2 = "d88e27aa-c46e-4436-a424-97a4e33712ca"
$hash1 = "5b66312beda68ab527bb8dfe2db90c72460d7a87"
$hash2 = "a93f81da4d30c56ad539b14ebefe0b3fdf0a4f15"

$status = "Ready for handoff"
$dateTimeFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(4,1).Value = "$file1.md"
$wsOverview.Cells.Item(4,2).Value = $status
$wsOverview.Cells.Item(4,3).Value = $status
$wsOverview.Cells.Item(4,4).Value = "2016-03-19 03:18:31"
$wsOverview.Cells.Item(4,4).NumberFormat = $dateTimeFormat

$wsOverview.Cells.Item(5,1).Value = "$file2.md"
$wsOverview.Cells.Item(5,2).Value = $status
$wsOverview.Cells.Item(5,3).Value = $status
$wsOverview.Cells.Item(5,4).Value = "2016-03-19 03:18:31"
$wsOverview.Cells.Item(5,4).NumberFormat = $dateTimeFormat

$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(4,1), "$mdBase$file1.md", "", "", "$file1.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(5,1), "$mdBase$file2.md", "", "", "$file2.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Cells.Item(4,1).Value = "$file1.md"
$wsZhCn.Cells.Item(4,2).Value = ".md"
$wsZhCn.Cells.Item(4,3).Value = $status
$wsZhCn.Cells.Item(4,4).Value = "$file1.$hash1.zh-cn.xlf"
$wsZhCn.Cells.Item(4,5).Value = "2016-03-19 03:18:23"
$wsZhCn.Cells.Item(4,5).NumberFormat = $dateTimeFormat
$wsZhCn.Cells.Item(4,8).Value = "0001-01-01 00:00:00"
$wsZhCn.Cells.Item(4,8).NumberFormat = $dateTimeFormat
$wsZhCn.Cells.Item(4,9).Value = ""
$wsZhCn.Cells.Item(4,10).Value = "Include"

$wsZhCn.Cells.Item(5,1).Value = "$file2.md"
$wsZhCn.Cells.Item(5,2).Value = ".md"
$wsZhCn.Cells.Item(5,3).Value = $status
$wsZhCn.Cells.Item(5,4).Value = "$file2.$hash2.zh-cn.xlf"
$wsZhCn.Cells.Item(5,5).Value = "2016-03-19 03:18:23"
$wsZhCn.Cells.Item(5,5).NumberFormat = $dateTimeFormat
$wsZhCn.Cells.Item(5,8).Value = "0001-01-01 00:00:00"
$wsZhCn.Cells.Item(5,8).NumberFormat = $dateTimeFormat
$wsZhCn.Cells.Item(5,9).Value = ""
$wsZhCn.Cells.Item(5,10).Value = "Include"

$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(4,1), "$mdBase$file1.md", "", "", "$file1.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(4,4), "$zhcnBase$file1.$hash1.zh-cn.xlf", "", "", "$file1.$hash1.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(5,1), "$mdBase$file2.md", "", "", "$file2.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(5,4), "$zhcnBase$file2.$hash2.zh-cn.xlf", "", "", "$file2.$hash2.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Cells.Item(4,1).Value = "$file1.md"
$wsDeDe.Cells.Item(4,2).Value = ".md"
$wsDeDe.Cells.Item(4,3).Value = $status
$wsDeDe.Cells.Item(4,4).Value = "$file1.$hash1.de-de.xlf"
$wsDeDe.Cells.Item(4,5).Value = "2016-03-19 03:18:31"
$wsDeDe.Cells.Item(4,5).NumberFormat = $dateTimeFormat
$wsDeDe.Cells.Item(4,8).Value = "0001-01-01 00:00:00"
$wsDeDe.Cells.Item(4,8).NumberFormat = $dateTimeFormat
$wsDeDe.Cells.Item(4,9).Value = ""
$wsDeDe.Cells.Item(4,10).Value = "Include"

$wsDeDe.Cells.Item(5,1).Value = "$file2.md"
$wsDeDe.Cells.Item(5,2).Value = ".md"
$wsDeDe.Cells.Item(5,3).Value = $status
$wsDeDe.Cells.Item(5,4).Value = "$file2.$hash2.de-de.xlf"
$wsDeDe.Cells.Item(5,5).Value = "2016-03-19 03:18:31"
$wsDeDe.Cells.Item(5,5).NumberFormat = $dateTimeFormat
$wsDeDe.Cells.Item(5,8).Value = "0001-01-01 00:00:00"
$wsDeDe.Cells.Item(5,8).NumberFormat = $dateTimeFormat
$wsDeDe.Cells.Item(5,9).Value = ""
$wsDeDe.Cells.Item(5,10).Value = "Include"

$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(4,1), "$mdBase$file1.md", "", "", "$file1.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(4,4), "$dedeBase$file1.$hash1.de-de.xlf", "", "", "$file1.$hash1.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(5,1), "$mdBase$file2.md", "", "", "$file2.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(5,4), "$dedeBase$file2.$hash2.de-de.xlf", "", "", "$file2.$hash2.de-de.xlf") | Out-Null

Write-Host "Report generated for handoff."
